$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '61.222.57'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '2.377.11'
$ws.Range("E3").Value = '  -4.20%  '
$ws.Range("E4").Value = '  +0.08%  '
Set-TextValue $ws.Range("D5") '549.03'
$ws.Range("E5").Value = '  -1.34%  '
Set-TextValue $ws.Range("D6") '141.63'
$ws.Range("E6").Value = '  -4.36%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -11.81%  '
$ws.Range("D9").Value = '2.375.94'
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("E10").Value = '  -2.71%  '
$ws.Range("E11").Value = '  +0.31%  '
Set-TextValue $ws.Range("D12") '5.30'
$ws.Range("E12").Value = '  -3.70%  '
Set-TextValue $ws.Range("D13") '0.346'
$ws.Range("E13").Value = '  -3.53%  '
Set-TextValue $ws.Range("D14") '25.31'
$ws.Range("E14").Value = '  -4.40%  '
$ws.Range("D15").Value = '2.808.13'
$ws.Range("E15").Value = '  -4.01%  '
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").Value = '61.168.19'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '2.380.94'
$ws.Range("E18").Value = '  -4.05%  '
Set-TextValue $ws.Range("D19") '10.71'
$ws.Range("E19").Value = '  -4.82%  '
$ws.Range("E20").Value = '  -2.88%  '
Set-TextValue $ws.Range("D21") '318.15'
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("E22").Value = '  -7.04%  '
$ws.Range("E23").Value = '  -0.13%  '
Set-TextValue $ws.Range("D24") '1.89'
$ws.Range("E24").Value = '  -0.40%  '
Set-TextValue $ws.Range("D25") '63.78'
$ws.Range("E25").Value = '  -1.09%  '
Set-TextValue $ws.Range("D26") '8.18'
$ws.Range("E26").Value = '  +4.25%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '2.496.40'
$ws.Range("E28").Value = '  -4.19%  '
$ws.Range("E29").Value = '  -9.36%  '
Set-TextValue $ws.Range("D30") '526.91'
$ws.Range("E30").Value = '  -7.04%  '
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("E32").Value = '  -6.10%  '
$ws.Range("E33").Value = '  -3.54%  '
$ws.Range("E34").Value = '  -5.41%  '
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  -5.81%  '
Set-TextValue $ws.Range("D38") '5.44'
$ws.Range("E38").Value = '  -8.87%  '
$ws.Range("E39").Value = '  +4.76%  '
Set-TextValue $ws.Range("D40") '0.374'
$ws.Range("E40").Value = '  -2.73%  '
Set-TextValue $ws.Range("D41") '18.12'
$ws.Range("E41").Value = '  -2.78%  '
Set-TextValue $ws.Range("D42") '139.94'
$ws.Range("E42").Value = '  -4.43%  '
$ws.Range("E43").Value = '  +0.03%  '
Set-TextValue $ws.Range("D44") '40.49'
Set-TextValue $ws.Range("D45") '2.14'
$ws.Range("E45").Value = '  -12.72%  '
Set-TextValue $ws.Range("D46") '140.43'
$ws.Range("E46").Value = '  -5.77%  '
$ws.Range("E47").Value = '  -1.63%  '
Set-TextValue $ws.Range("D48") '20.08'
$ws.Range("E48").Value = '  -9.21%  '
$ws.Range("E49").Value = '  -4.78%  '
Set-TextValue $ws.Range("D50") '0.574'
$ws.Range("E50").Value = '  -4.11%  '
Set-TextValue $ws.Range("D51") '0.0904'
$ws.Range("E51").Value = '  -4.45%  '
